# Added filtering options for the Component Analysis
# Trim the trailing forecast-horizon columns from rows 2, 3, 5, 6 and 7 so
# that every row follows the same staircase pattern as the rest of the
# error-series table (row 2 -> through column I, row 3 -> through column H,
# row 5 -> through column J, row 6 -> through column I, row 7 -> through column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5:K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
